$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1): update "想去人数" (want-to-go count) values in column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 129
$ws1.Range("F4").Value = 189
$ws1.Range("F5").Value = 3429
$ws1.Range("F6").Value = 348
$ws1.Range("F7").Value = 19
$ws1.Range("F8").Value = 425

# Sheet "全部类型" (index 4): same events appear here, rows shifted because
# of the extra rows coming from the "演出" sheet, so F7/F8 differ but
# F3-F6 and F9/F10 correspond to the same events as above.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 129
$ws4.Range("F4").Value = 189
$ws4.Range("F5").Value = 3429
$ws4.Range("F6").Value = 348
$ws4.Range("F9").Value = 19
$ws4.Range("F10").Value = 425
